$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 264 (A1:T264). Append three new weekly
# price rows (265-267) for "Damasco" / variety "Dina" at the same market,
# date 2023-01-13 (serial 44939), mirroring the structure of the existing
# data rows.

$newRows = @(
    @{ Row = 265; L = "Especial"; M = 155; N = 17000; O = 17000; P = 17000; S = 1062 },
    @{ Row = 266; L = "Primera";  M = 275; N = 15000; O = 15000; P = 15000; S = 938  },
    @{ Row = 267; L = "Segunda";  M = 250; N = 12000; O = 12000; P = 12000; S = 750  }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $ws.Range("A$r").Value = 6
    $ws.Range("B$r").Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Range("C$r").Value = "Metropolitana"

    $ws.Range("D$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("D$r").Value = 44939

    $ws.Range("E$r").Value = 13
    $ws.Range("F$r").Value = "Fruta"
    $ws.Range("G$r").Value = 100103
    $ws.Range("H$r").Value = "Frutos de hueso (carozo)"
    $ws.Range("I$r").Value = 100103003
    $ws.Range("J$r").Value = "Damasco"
    $ws.Range("K$r").Value = "Dina"
    $ws.Range("L$r").Value = $item.L
    $ws.Range("M$r").Value = $item.M
    $ws.Range("N$r").Value = $item.N
    $ws.Range("O$r").Value = $item.O
    $ws.Range("P$r").Value = $item.P
    $ws.Range("Q$r").Value = "`$/caja 16 kilos"
    $ws.Range("R$r").Value = "Región Metropolitana"
    $ws.Range("S$r").Value = $item.S
    $ws.Range("T$r").Value = 16
}

Write-Host "Appended rows 265-267"
